$wb = $excel.ActiveWorkbook

# Remove the "Professors" sheet - only "Courses" remains.
[void]$wb.Worksheets("Professors").Delete()

$ws = $wb.Worksheets("Courses")

# --- Header row: drop the old "Course Number" column, shifting everything
#     one column to the left, and drop the "(to be sent in private)" suffixes. ---
$ws.Range("A1:J1").Value = $null
$ws.Range("A1").Value = "Course Name"
$ws.Range("B1").Value = "Whatsapp Group URL"
$ws.Range("C1").Value = "Faculty"
$ws.Range("D1").Value = "Telegram Group URL"
$ws.Range("E1").Value = "Drive Link"
$ws.Range("F1").Value = "Course info STKFUPM URL"
$ws.Range("G1").Value = "Course overall info"
$ws.Range("H1").Value = "Study Tips"

# --- Row 2: ICS104 sample data ---
$ws.Range("A2").Value = "ICS104"
$ws.Range("B2").Value = "http://resources.stkfupm.com/?dir=uploads/ICS/ICS104"
$ws.Range("C2").Value = "Dr.1/Dr.2/Dr.3/Dr.4/Dr.5/Dr.6/Dr.7/Dr.8"
$ws.Range("D2").Value = "http://www.kfupm.edu.sa/departments/ics/Pages/en/Undergraduate-Courses.aspx"
$ws.Range("E2").Value = "https://drive.google.com/drive/folders/1_KyrFSdSmgG-SyqpNFuwXKB_PKQlJ97X"
$ws.Range("F2").Value = "http://www.kfupm.edu.sa/departments/ics/Pages/en/Undergraduate-Courses.aspx"
$ws.Range("G2").Value = "Info About ICS104"
$ws.Range("H2").Value = "Study Tips About ICS104"

# --- Row 3: EE311 sample data ---
$ws.Range("A3").Value = "EE311"
$ws.Range("B3").Value = "https://www.stkfupm.com/forum/threads/ee-311-fundamentals-of-electrical-engineering-design.110154/"
$ws.Range("C3").Value = "Dr.3"
$ws.Range("D3").Value = "https://bulletin.kfupm.edu.sa/course-detail?course_code=EE311"
$ws.Range("E3").Value = "https://drive.google.com/drive/folders/1EzdZ5iWV9rciGxuicqIcGhXEy2tEM2iA"
$ws.Range("F3").Value = "https://bulletin.kfupm.edu.sa/course-detail?course_code=EE311"
$ws.Range("G3").Value = "Info About EE311"
$ws.Range("H3").Value = "Study Tips About EE311"
